# The document carries three inline logo pictures (two Pearson logos in
# the footers, one BTEC logo in the header). This change only swaps the
# internal drawing "name" label Word stores for each picture - the image
# bytes, positions, alt text, etc. are untouched:
#
#   - both Pearson logo pictures (in the footers): "image2.png" -> "image1.png"
#   - the BTEC logo picture (in the header):       "image1.jpg" -> "image2.jpg"
#
# InlineShape itself exposes no settable Name property (this matches real
# Word - only floating Shapes have .Name), so each picture is converted to
# a floating shape long enough to rename it, then converted straight back
# to an inline picture in place.

function Rename-LogoPicture($range, $targetAltText, $newName) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $inlineShape = $shapes.Item($i)
        if ($inlineShape.AlternativeText -eq $targetAltText) {
            $floatingShape = $inlineShape.ConvertToShape()
            $floatingShape.Name = $newName
            $floatingShape.ConvertToInlineShape() | Out-Null
        }
    }
}

$d = $word.ActiveDocument
$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecAlt = "BTec_Logo-Orange"

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($fi = 1; $fi -le 3; $fi++) {
        $footer = $sec.Footers.Item($fi)
        if ($footer.Exists) {
            Rename-LogoPicture $footer.Range $pearsonAlt "image1.png"
        }
    }

    for ($hi = 1; $hi -le 3; $hi++) {
        $header = $sec.Headers.Item($hi)
        if ($header.Exists) {
            Rename-LogoPicture $header.Range $btecAlt "image2.jpg"
        }
    }
}
